$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1206.8334
$ws.Range("I32").Value = 948
$ws.Range("J32").Value = 1336.25
$ws.Range("K32").Value = 948
$ws.Range("L32").Value = 1336.25
$ws.Range("M32").Value = -622
$ws.Range("N32").Value = -1988.25
$ws.Range("H137").Value = 40175.04
$ws.Range("I137").Value = 126321.75
$ws.Range("J137").Value = 1887.6111
$ws.Range("K137").Value = 378965.25
$ws.Range("L137").Value = 5662.8333
$ws.Range("M137").Value = -376415.25
$ws.Range("N137").Value = -10762.8333
$ws.Range("H138").Value = 3898.5103
$ws.Range("I138").Value = 948.6875
$ws.Range("J138").Value = 4474.0854
$ws.Range("K138").Value = 2846.0625
$ws.Range("L138").Value = 13422.2562
$ws.Range("M138").Value = 2293.9375
$ws.Range("N138").Value = -23702.2562
$ws.Range("H141").Value = 2279.0293
$ws.Range("I141").Value = 1778.2084
$ws.Range("J141").Value = 3481
$ws.Range("K141").Value = 5334.6252
$ws.Range("L141").Value = 10443
$ws.Range("M141").Value = -154.6252000000004
$ws.Range("N141").Value = -20803

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 3537.4
$ws.Range("I5").Value = 4299.25
$ws.Range("K5").Value = 4299.25
$ws.Range("M5").Value = -4187.25
$ws.Range("H32").Value = 1823219
$ws.Range("I32").Value = 1933161.9
$ws.Range("J32").Value = 256532
$ws.Range("K32").Value = 1933161.9
$ws.Range("L32").Value = 256532
$ws.Range("M32").Value = -1932874.9
$ws.Range("N32").Value = -257106
$ws.Range("H35").Value = 33123.8
$ws.Range("I35").Value = 1845.6666
$ws.Range("K35").Value = 1845.6666
$ws.Range("M35").Value = -1439.6666
$ws.Range("H62").Value = 90244.5
$ws.Range("J62").Value = 90244.5
$ws.Range("L62").Value = 90244.5
$ws.Range("N62").Value = -91492.5
$ws.Range("H65").Value = 90244.5
$ws.Range("J65").Value = 90244.5
$ws.Range("L65").Value = 270733.5
$ws.Range("N65").Value = -276973.5
$ws.Range("H75").Value = 500045100
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 500045100
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H81").Value = 336356000
$ws.Range("I81").Value = 9000000
$ws.Range("J81").Value = 500034000
$ws.Range("K81").Value = 9000000
$ws.Range("L81").Value = 500034000
$ws.Range("M81").Value = -8999002
$ws.Range("N81").Value = -500035996
$ws.Range("H84").Value = 336356000
$ws.Range("I84").Value = 9000000
$ws.Range("J84").Value = 500034000
$ws.Range("K84").Value = 27000000
$ws.Range("L84").Value = 1500102000
$ws.Range("M84").Value = -26995008
$ws.Range("N84").Value = -1500111984
$ws.Range("H87").Value = 40175
$ws.Range("J87").Value = 40175
$ws.Range("L87").Value = 40175
$ws.Range("N87").Value = -42671
$ws.Range("H90").Value = 40175
$ws.Range("J90").Value = 40175
$ws.Range("L90").Value = 120525
$ws.Range("N90").Value = -133005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 3537.4
$ws.Range("I4").Value = 4299.25
$ws.Range("K4").Value = 4299.25
$ws.Range("M4").Value = -4184.25
$ws.Range("H22").Value = 487.94446
$ws.Range("I22").Value = 427
$ws.Range("K22").Value = 427
$ws.Range("M22").Value = -254
$ws.Range("H37").Value = 40000
$ws.Range("J37").Value = 40000
$ws.Range("L37").Value = 40000
$ws.Range("N37").Value = -40274
$ws.Range("H134").Value = 33402646
$ws.Range("I134").Value = 55556180
$ws.Range("K134").Value = 166668540
$ws.Range("M134").Value = -166666005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1723
$ws.Range("I7").Value = 2534.5
$ws.Range("K7").Value = 2534.5
$ws.Range("M7").Value = -2421.5
$ws.Range("H31").Value = 43973.4
$ws.Range("I31").Value = 144633.58
$ws.Range("J31").Value = 4827.778
$ws.Range("K31").Value = 144633.58
$ws.Range("L31").Value = 4827.778
$ws.Range("M31").Value = -144338.58
$ws.Range("N31").Value = -5417.778
$ws.Range("H34").Value = 43973.4
$ws.Range("I34").Value = 144633.58
$ws.Range("J34").Value = 4827.778
$ws.Range("K34").Value = 144633.58
$ws.Range("L34").Value = 4827.778
$ws.Range("M34").Value = -144431.58
$ws.Range("N34").Value = -5231.778

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5883081
$ws.Range("I4").Value = 25000170
$ws.Range("J4").Value = 899.7692
$ws.Range("K4").Value = 75000510
$ws.Range("L4").Value = 2699.3076
$ws.Range("M4").Value = -75000398
$ws.Range("N4").Value = -2923.3076
$ws.Range("H5").Value = 8934.416999999999
$ws.Range("I5").Value = 580.8889
$ws.Range("K5").Value = 1742.6667
$ws.Range("M5").Value = -1630.6667
$ws.Range("H135").Value = 8934.416999999999
$ws.Range("I135").Value = 580.8889
$ws.Range("K135").Value = 5228.0001
$ws.Range("M135").Value = -2693.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 10767.125
$ws.Range("I43").Value = 1000
$ws.Range("K43").Value = 1000
$ws.Range("M43").Value = -849
$ws.Range("H88").Value = 59900
$ws.Range("J88").Value = 59900
$ws.Range("L88").Value = 59900
$ws.Range("N88").Value = -60802
$ws.Range("H91").Value = 59900
$ws.Range("J91").Value = 59900
$ws.Range("L91").Value = 59900
$ws.Range("N91").Value = -63020
$ws.Range("H132").Value = 24118.564
$ws.Range("I132").Value = 2242.4644
$ws.Range("J132").Value = 58148.055
$ws.Range("K132").Value = 6727.3932
$ws.Range("L132").Value = 174444.165
$ws.Range("M132").Value = -4197.3932
$ws.Range("N132").Value = -179504.165

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 567.25
$ws.Range("J22").Value = 632
$ws.Range("L22").Value = 632
$ws.Range("N22").Value = -1222
$ws.Range("H27").Value = 567.25
$ws.Range("J27").Value = 632
$ws.Range("L27").Value = 632
$ws.Range("N27").Value = -846
$ws.Range("H46").Value = 967.7692
$ws.Range("I46").Value = 798.2
$ws.Range("J46").Value = 1073.75
$ws.Range("K46").Value = 798.2
$ws.Range("L46").Value = 1073.75
$ws.Range("M46").Value = -610.2
$ws.Range("N46").Value = -1449.75
$ws.Range("H61").Value = 4112.5
$ws.Range("I61").Value = 5950
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 5950
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -5748
$ws.Range("N61").Value = -3904
$ws.Range("H113").Value = 4112.5
$ws.Range("I113").Value = 5950
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 5950
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -3780
$ws.Range("N113").Value = -7840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 45406
$ws.Range("J64").Value = 45406
$ws.Range("L64").Value = 45406
$ws.Range("N64").Value = -45902
$ws.Range("H67").Value = 45406
$ws.Range("J67").Value = 45406
$ws.Range("L67").Value = 45406
$ws.Range("N67").Value = -47122
